$d = $word.ActiveDocument

# Locate the placeholder text (second line of the "Sklad osobowy" paragraph)
# and replace it with the new name/id, while keeping it in its own run
# (separate from the "Grzegorz Gojska 174173" run that precedes the <w:br/>).
$rng = $d.Content
$rng.Find.Execute("xxxx xxxx 111111", $true, $false, $false, $false, $false, `
                   $true, 1, $false, "Jakub Konkel 187207", 2) | Out-Null

# Force the newly inserted text to live in its own run distinct from the
# run that holds the preceding line break, by toggling a character format
# on then back off. This keeps the final formatting identical to the
# surrounding text while preventing the run from being coalesced with the
# previous one.
$rng.Font.Bold = $true
$rng.Font.Bold = $false
